$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.230" or
# "12.70" keep their exact literal representation instead of being
# auto-converted to numbers (which would drop trailing zeros / use
# scientific notation for very small numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.683.41"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.655.04"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "303.02"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.3835"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "0.3606"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "51.11"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "0.08194"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "1.230"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "6.453"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "7.440"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "0.00001221"
$ws.Range("D17").Value = "1.656.28"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "97.61"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "0.07036"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "6.781"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").Value = "17.55"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "12.70"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").Value = "23.683.54"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "2.487"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "3.030"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "21.26"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "152.88"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "5.240"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "134.02"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "1.839.86"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "7.089"
$ws.Range("E32").Value = "  +8.89%  "
$ws.Range("D33").Value = "2.256"
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").Value = "12.11"
$ws.Range("E34").Value = "  +5.74%  "
$ws.Range("D35").Value = "1.059"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").Value = "0.02806"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "0.2503"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "0.08810"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "6.078"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "0.06987"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("D42").Value = "0.6988"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").Value = "1.339"
$ws.Range("D44").Value = "15.93"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").Value = "0.6504"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "2.299"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "3.964"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "0.07894"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "128.20"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "1.178"
$ws.Range("E51").Value = "  -0.87%  "
